# Update cryptocurrency price/volume data (cryptos list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.539.52'
$ws.Range('E2').Value = '  -0.10%  '
$ws.Range('D3').Value = '1.824.58'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('D4').Formula = "'1.001"
$ws.Range('D5').Formula = "'315.52"
$ws.Range('E5').Value = '  -0.64%  '
$ws.Range('E6').Value = '  +0.06%  '
$ws.Range('D7').Formula = "'0.5112"
$ws.Range('E7').Value = '  -5.89%  '
$ws.Range('D8').Formula = "'0.3959"
$ws.Range('E8').Value = '  -1.89%  '
$ws.Range('D9').Formula = "'0.08228"
$ws.Range('E9').Value = '  +7.16%  '
$ws.Range('D10').Formula = "'1.113"
$ws.Range('E10').Value = '  -0.75%  '
$ws.Range('E11').Value = '  -0.44%  '
$ws.Range('D12').Formula = "'21.20"
$ws.Range('D13').Formula = "'6.334"
$ws.Range('E13').Value = '  +0.12%  '
$ws.Range('E14').Value = '  +0.09%  '
$ws.Range('D15').Formula = "'7.534"
$ws.Range('E15').Value = '  -1.22%  '
$ws.Range('D16').Value = '1.822.07'
$ws.Range('E16').Value = '  -0.24%  '
$ws.Range('D17').Formula = "'0.00001129"
$ws.Range('E17').Value = '  +3.40%  '
$ws.Range('D18').Formula = "'92.71"
$ws.Range('E18').Value = '  +3.01%  '
$ws.Range('D19').Formula = "'0.06666"
$ws.Range('E19').Value = '  +0.86%  '
$ws.Range('D20').Formula = "'17.84"
$ws.Range('E20').Value = '  +0.25%  '
$ws.Range('E21').Value = '  +0.01%  '
$ws.Range('D22').Formula = "'6.101"
$ws.Range('E22').Value = '  +0.54%  '
$ws.Range('D23').Value = '28.570.74'
$ws.Range('E23').Value = '  +0.02%  '
$ws.Range('D24').Formula = "'11.41"
$ws.Range('E24').Value = '  +2.13%  '
$ws.Range('E25').Value = '  -0.60%  '
$ws.Range('D26').Formula = "'21.50"
$ws.Range('E26').Value = '  +3.44%  '
$ws.Range('D27').Formula = "'156.79"
$ws.Range('E27').Value = '  -0.77%  '
$ws.Range('D28').Value = '2.032.71'
$ws.Range('E28').Value = '  -0.18%  '
$ws.Range('D29').Formula = "'2.405"
$ws.Range('E29').Value = '  -2.15%  '
$ws.Range('D30').Formula = "'126.81"
$ws.Range('E30').Value = '  +2.19%  '
$ws.Range('D31').Formula = "'1.113"
$ws.Range('E31').Value = '  -1.53%  '
$ws.Range('D32').Formula = "'0.1093"
$ws.Range('E32').Value = '  -1.43%  '
$ws.Range('D33').Formula = "'5.767"
$ws.Range('E33').Value = '  +1.38%  '
$ws.Range('D34').Formula = "'3.661"
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('D35').Formula = "'0.07070"
$ws.Range('E35').Value = '  -4.23%  '
$ws.Range('D36').Formula = "'0.2234"
$ws.Range('E36').Value = '  -0.60%  '
$ws.Range('D37').Formula = "'5.293"
$ws.Range('E37').Value = '  +1.51%  '
$ws.Range('D38').Formula = "'0.02361"
$ws.Range('E38').Value = '  +0.13%  '
$ws.Range('D39').Formula = "'8.848"
$ws.Range('E39').Value = '  -0.60%  '
$ws.Range('D40').Formula = "'0.6327"
$ws.Range('E40').Value = '  +0.40%  '
$ws.Range('D41').Formula = "'11.32"
$ws.Range('E41').Value = '  -0.41%  '
$ws.Range('D42').Formula = "'1.184"
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('E43').Value = '  +0.05%  '
$ws.Range('E44').Value = '  -0.27%  '
$ws.Range('D45').Formula = "'13.58"
$ws.Range('E45').Value = '  +1.19%  '
$ws.Range('D46').Formula = "'0.5950"
$ws.Range('E46').Value = '  +1.08%  '
$ws.Range('D47').Formula = "'3.735"
$ws.Range('E47').Value = '  +0.70%  '
$ws.Range('D48').Formula = "'125.69"
$ws.Range('E48').Value = '  +0.41%  '
$ws.Range('D49').Formula = "'1.997"
$ws.Range('E49').Value = '  -0.33%  '
$ws.Range('D50').Formula = "'1.192"
$ws.Range('E50').Value = '  -0.54%  '
$ws.Range('D51').Formula = "'0.06920"
$ws.Range('E51').Value = '  +0.14%  '
